# Convert HOUR_APPR_PROCESS_START (column V) values from numbers to
# strings formatted like "H:00:00" (a time-like text value), e.g. 5 -> "5:00:00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on the sheet's dimension.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("V$row")
    $hour = $cell.Value2
    if ($null -ne $hour -and $hour -ne "") {
        $cell.Value = "$($hour):00:00"
    }
}
